# Apply the changes described in the commit:
#  - Steps!E3:E52 become an incrementing sequence (1..50) computed with
#    formulas instead of hard-coded decimal "major.minor" step numbers.
#  - Steps!I52 now looks up the second guest (guest2) instead of the
#    first guest (guest1).
#  - The Steps sheet tab becomes the active/selected sheet (instead of
#    Speech), with the view scrolled down and the selection moved to I53.

$wb = $excel.ActiveWorkbook
$wsSteps = $wb.Worksheets.Item("Steps")

# --- E column: replace literal step numbers with a running count ---
# E3 = E2 + 1
$wsSteps.Range("E3").Formula = "=E2+1"

# E4:E52 = previous cell + 1 (Excel will store this as one shared formula,
# matching the "t=shared" group used for the rest of the sheet)
$wsSteps.Range("E4:E52").Formula = "=E3+1"

# --- I52: point at the second guest instead of the first guest ---
$wsSteps.Range("I52").Formula = "=guest2"

# --- Selection / active sheet housekeeping ---
# Make the Steps sheet the active tab (it was Speech before), scroll it
# down so row 22 is at the top, and leave the selection on I53.
$wsSteps.Activate()
$wn = $excel.ActiveWindow
$wn.ScrollRow = 22
$wn.ScrollColumn = 1
$wsSteps.Range("I53").Select()

$excel.CalculateFull()
